# "Generate Report for Handoff" — the localization-status report was
# regenerated: the translation status flips from "In Translation" to
# "Ready for handoff" and the handoff/generation timestamps advance by
# under a minute. Excel's real AutoFit then widens the (now slightly
# longer) status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"          # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"          # de-de status
$overview.Range("G2").Value = "2016-08-12 22:49:23"        # Latest HO Xliff Generate Date

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"               # Status
$zhcn.Range("H2").Value = "2016-08-12 22:49:17"             # Latest Handoff Datetime

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"                # Status
$dede.Range("H2").Value = "2016-08-12 22:49:23"              # Latest Handoff Datetime

# --- Re-fit the Status columns now that their text grew longer --------
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
